$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("A2").Value = "12 Nov 2025, 11:56 AM"

$ws = $wb.Worksheets.Item("2 Week Return")
$ws.Range("C2").Value = 66.7595
$ws.Range("C3").Value = 55.525
$ws.Range("C4").Value = 44.5092
$ws.Range("C5").Value = 44.4151
$ws.Range("C6").Value = 42.7197
$ws.Range("C7").Value = 40.6424
$ws.Range("C8").Value = 36.0648
$ws.Range("B9").Value = "DREDGECORP"
$ws.Range("C9").Value = 35.551
$ws.Range("B10").Value = "MAHASTEEL"
$ws.Range("C10").Value = 35.3849
$ws.Range("C11").Value = 33.4403
$ws.Range("C12").Value = 31.7995
$ws.Range("C13").Value = 26.3847
$ws.Range("C14").Value = 26.3067
$ws.Range("C15").Value = 25.7344
$ws.Range("C18").Value = 24.2331
$ws.Range("C19").Value = 23.8498
$ws.Range("C20").Value = 23.4318
$ws.Range("C21").Value = 23.4168
$ws.Range("C22").Value = 22.9423
$ws.Range("C24").Value = 22.1654
$ws.Range("C25").Value = 22.0333
$ws.Range("C27").Value = 21.2152
$ws.Range("C29").Value = 20.6548
$ws.Range("C30").Value = 20.5928
$ws.Range("C32").Value = 20.3118
$ws.Range("C33").Value = 19.7975
$ws.Range("B35").Value = "SHRIRAMFIN"
$ws.Range("C35").Value = 19.0795
$ws.Range("B36").Value = "MTARTECH"
$ws.Range("C36").Value = 19.0456
$ws.Range("C39").Value = 18.4943
$ws.Range("C40").Value = 18.3235
$ws.Range("C41").Value = 18.1663
$ws.Range("B42").Value = "KAPSTON"
$ws.Range("C42").Value = 18.0357
$ws.Range("B43").Value = "PARAGMILK"
$ws.Range("C43").Value = 17.9326
$ws.Range("C44").Value = 17.4752
$ws.Range("C45").Value = 16.9214
$ws.Range("C46").Value = 16.9097
$ws.Range("C48").Value = 16.6632
$ws.Range("B49").Value = "INDIGOPNTS"
$ws.Range("C49").Value = 16.2753
$ws.Range("B50").Value = "GSLSU"
$ws.Range("C50").Value = 16.2581
$ws.Range("B51").Value = "SILVERTUC"
$ws.Range("C51").Value = 16.1563
$ws.Range("B52").Value = "VOLTAMP"
$ws.Range("C52").Value = 16.0672
$ws.Range("B53").Value = "SCI"
$ws.Range("C53").Value = 15.9774
$ws.Range("B54").Value = "BORANA"
$ws.Range("C54").Value = 15.9442
$ws.Range("B55").Value = "NAVINFLUOR"
$ws.Range("C55").Value = 15.8874
$ws.Range("C56").Value = 15.6617
$ws.Range("B57").Value = "STAR"
$ws.Range("C57").Value = 15.4904
$ws.Range("B58").Value = "CARYSIL"
$ws.Range("C58").Value = 15.4824
$ws.Range("C60").Value = 15.0394
$ws.Range("C61").Value = 14.6863
$ws.Range("C62").Value = 14.5395
$ws.Range("C63").Value = 14.2822
$ws.Range("C65").Value = 13.9092
$ws.Range("B67").Value = "INTELLECT"
$ws.Range("C67").Value = 13.3904
$ws.Range("B68").Value = "SHREEJISPG"
$ws.Range("C68").Value = 13.3731
$ws.Range("B69").Value = "FSL"
$ws.Range("C69").Value = 13.2691
$ws.Range("B70").Value = "CUPID"
$ws.Range("C70").Value = 13.2155
$ws.Range("C71").Value = 13.1565
$ws.Range("B72").Value = "IVALUE"
$ws.Range("C72").Value = 13.0996
$ws.Range("B73").Value = "INDUSTOWER"
$ws.Range("C73").Value = 12.9614
$ws.Range("B74").Value = "CANBK"
$ws.Range("C74").Value = 12.9484
$ws.Range("C75").Value = 12.6544
$ws.Range("B76").Value = "CUB"
$ws.Range("C76").Value = 12.6072

$ws = $wb.Worksheets.Item("Industry Analysis")
$ws.Range("C2").Value = 0.4615
$ws.Range("D2").Value = 20.3167
$ws.Range("E2").Value = 20.8221
$ws.Range("F2").Value = -13.8333
$ws.Range("C4").Value = 2.9066
$ws.Range("D4").Value = 14.7975
$ws.Range("E4").Value = 17.989
$ws.Range("C5").Value = 14.2162
$ws.Range("D5").Value = 11.9718
$ws.Range("E5").Value = 7.9908
$ws.Range("F5").Value = 5.1258
$ws.Range("C6").Value = 3.6392
$ws.Range("D6").Value = 11.0343
$ws.Range("E6").Value = 22.2652
$ws.Range("F6").Value = 28.3195
$ws.Range("B8").Value = "electronics - components"
$ws.Range("C8").Value = 0.9685
$ws.Range("D8").Value = 6.3909
$ws.Range("E8").Value = 1.4359
$ws.Range("F8").Value = 12.6563
$ws.Range("G8").Value = 32.7477
$ws.Range("H8").Value = 69.8873
$ws.Range("I8").Value = 40.947
$ws.Range("J8").Value = 59.6719
$ws.Range("K8").Value = 55.4676
$ws.Range("B9").Value = "moulded luggage"
$ws.Range("C9").Value = -0.4022
$ws.Range("D9").Value = 6.3705
$ws.Range("E9").Value = 2.0704
$ws.Range("F9").Value = -16.2487
$ws.Range("G9").Value = 33.3308
$ws.Range("H9").Value = 131.9308
$ws.Range("I9").Value = 72.0789
$ws.Range("J9").Value = 61.9537
$ws.Range("K9").Value = 54.3963
$ws.Range("B11").Value = "glass & glass products"
$ws.Range("C11").Value = -0.3289
$ws.Range("D11").Value = 6.0721
$ws.Range("E11").Value = 8.2971
$ws.Range("F11").Value = 24.3643
$ws.Range("G11").Value = 27.404
$ws.Range("H11").Value = -13.5162
$ws.Range("I11").Value = 107.7717
$ws.Range("J11").Value = 50.9508
$ws.Range("K11").Value = 22.7727
$ws.Range("B12").Value = "electric equipment"
$ws.Range("C12").Value = 0.9518
$ws.Range("D12").Value = 5.747
$ws.Range("E12").Value = 24.0972
$ws.Range("F12").Value = 53.0786
$ws.Range("G12").Value = 173.7698
$ws.Range("H12").Value = 57.2692
$ws.Range("I12").Value = 95.4705
$ws.Range("J12").Value = 74.5879
$ws.Range("K12").Value = 28.2089
$ws.Range("B13").Value = "bearings"
$ws.Range("C13").Value = -0.5602
$ws.Range("D13").Value = 5.3803
$ws.Range("E13").Value = 3.7254
$ws.Range("F13").Value = 1.5374
$ws.Range("G13").Value = -16.7188
$ws.Range("H13").Value = 120.0131
$ws.Range("I13").Value = 58.3623
$ws.Range("J13").Value = 32.445
$ws.Range("K13").Value = 17.366
$ws.Range("C15").Value = 0.8073
$ws.Range("D15").Value = 4.3677
$ws.Range("E15").Value = 6.1852
$ws.Range("F15").Value = 6.4021
$ws.Range("C16").Value = -0.068
$ws.Range("D16").Value = 4.282
$ws.Range("E16").Value = 9.2103
$ws.Range("F16").Value = 52.7903
$ws.Range("C17").Value = 0.6312
$ws.Range("D17").Value = 4.2178
$ws.Range("E17").Value = 3.4489
$ws.Range("F17").Value = 27.7708
$ws.Range("C18").Value = 6.8805
$ws.Range("D18").Value = 3.6473
$ws.Range("E18").Value = 1.4714
$ws.Range("F18").Value = -54.3434
$ws.Range("C19").Value = 0.2187
$ws.Range("D19").Value = 3.3916
$ws.Range("E19").Value = 3.0356
$ws.Range("F19").Value = -24.5109
$ws.Range("B20").Value = "computers - software - large"
$ws.Range("C20").Value = 2.5758
$ws.Range("D20").Value = 2.937
$ws.Range("E20").Value = 2.1972
$ws.Range("F20").Value = -11.1475
$ws.Range("G20").Value = -9.848800000000001
$ws.Range("H20").Value = 42.6599
$ws.Range("I20").Value = 139.3123
$ws.Range("J20").Value = 39.4476
$ws.Range("K20").Value = 25.3058
$ws.Range("B21").Value = "ceramics - tiles / sanitaryware"
$ws.Range("C21").Value = 0.1778
$ws.Range("D21").Value = 2.8754
$ws.Range("E21").Value = 6.5721
$ws.Range("F21").Value = 0.61
$ws.Range("G21").Value = -10.7437
$ws.Range("H21").Value = 48.5207
$ws.Range("I21").Value = -50.3371
$ws.Range("J21").Value = -9.548999999999999
$ws.Range("K21").Value = 11.6661
$ws.Range("F22").Value = 12.8436
$ws.Range("C23").Value = 1.8077
$ws.Range("D23").Value = 2.8171
$ws.Range("E23").Value = 1.0727
$ws.Range("F23").Value = 25.1253
$ws.Range("B24").Value = "steel - medium / small"
$ws.Range("C24").Value = 0.4169
$ws.Range("D24").Value = 2.8145
$ws.Range("E24").Value = 0.0649
$ws.Range("F24").Value = 14.2264
$ws.Range("G24").Value = -19.3113
$ws.Range("H24").Value = 92.69580000000001
$ws.Range("I24").Value = -9.5832
$ws.Range("J24").Value = 22.3385
$ws.Range("K24").Value = 30.7144
$ws.Range("C26").Value = 2.5529
$ws.Range("D26").Value = 2.3186
$ws.Range("E26").Value = -1.5182
$ws.Range("F26").Value = -6.5656
$ws.Range("C27").Value = 0.6029
$ws.Range("D27").Value = 2.2628
$ws.Range("E27").Value = -1.9943
$ws.Range("F27").Value = -15.3393
$ws.Range("B28").Value = "auto ancillaries"
$ws.Range("C28").Value = 0.1099
$ws.Range("D28").Value = 2.2379
$ws.Range("E28").Value = 2.3122
$ws.Range("F28").Value = 3.9501
$ws.Range("G28").Value = 46.2142
$ws.Range("H28").Value = 36.882
$ws.Range("I28").Value = 67.149
$ws.Range("J28").Value = 28.2354
$ws.Range("K28").Value = 13.9865
$ws.Range("B29").Value = "oil drilling / allied services"
$ws.Range("C29").Value = 1.6745
$ws.Range("D29").Value = 2.202
$ws.Range("E29").Value = 2.2845
$ws.Range("F29").Value = 2.9019
$ws.Range("G29").Value = 73.6195
$ws.Range("H29").Value = 78.74639999999999
$ws.Range("I29").Value = 85.1955
$ws.Range("J29").Value = 42.4825
$ws.Range("K29").Value = 16.4707
$ws.Range("B30").Value = "cables - power"
$ws.Range("C30").Value = 0.6332
$ws.Range("D30").Value = 2.201
$ws.Range("E30").Value = 0.0519
$ws.Range("F30").Value = 5.9905
$ws.Range("G30").Value = 32.5898
$ws.Range("H30").Value = 113.493
$ws.Range("I30").Value = 137.884
$ws.Range("J30").Value = 58.5243
$ws.Range("K30").Value = 72.6962
$ws.Range("C31").Value = 0.5036
$ws.Range("D31").Value = 2.0453
$ws.Range("E31").Value = -1.6988
$ws.Range("F31").Value = -15.6692
$ws.Range("C32").Value = 0.0406
$ws.Range("D32").Value = 1.9804
$ws.Range("E32").Value = 3.9843
$ws.Range("F32").Value = 26.4643
$ws.Range("C33").Value = -0.045
$ws.Range("D33").Value = 1.9435
$ws.Range("E33").Value = -0.0169
$ws.Range("F33").Value = 1.0502
$ws.Range("B35").Value = "tea"
$ws.Range("C35").Value = 0.7873
$ws.Range("D35").Value = 1.7753
$ws.Range("E35").Value = 21.355
$ws.Range("F35").Value = 40.5848
$ws.Range("G35").Value = 15.2854
$ws.Range("H35").Value = 21.5655
$ws.Range("I35").Value = 60.3474
$ws.Range("J35").Value = 30.7101
$ws.Range("K35").Value = 18.4152
$ws.Range("B37").Value = "recreation / amusement parks"
$ws.Range("C37").Value = -1.2843
$ws.Range("D37").Value = 1.6782
$ws.Range("E37").Value = -0.2325
$ws.Range("F37").Value = -26.7425
$ws.Range("G37").Value = 18.4499
$ws.Range("H37").Value = 59.973
$ws.Range("I37").Value = 86.66670000000001
$ws.Range("J37").Value = 86.5406
$ws.Range("K37").Value = 26.5444
$ws.Range("C38").Value = 1.1756
$ws.Range("D38").Value = 1.4661
$ws.Range("E38").Value = 6.0519
$ws.Range("F38").Value = -49.1741
$ws.Range("C39").Value = 0.5689
$ws.Range("D39").Value = 1.3911
$ws.Range("E39").Value = 2.2269
$ws.Range("F39").Value = -13.4028
$ws.Range("C40").Value = -0.1969
$ws.Range("D40").Value = 1.1113
$ws.Range("E40").Value = -4.4304
$ws.Range("F40").Value = -27.6481
$ws.Range("C41").Value = -0.1918
$ws.Range("D41").Value = 0.8786
$ws.Range("E41").Value = -3.5277
$ws.Range("F41").Value = 43.8057
$ws.Range("C42").Value = 0.3536
$ws.Range("D42").Value = 0.5603
$ws.Range("E42").Value = -1.0299
$ws.Range("F42").Value = -18.1159
$ws.Range("B43").Value = "food - processing - indian"
$ws.Range("C43").Value = 1.8503
$ws.Range("D43").Value = 0.5254
$ws.Range("E43").Value = -6.6734
$ws.Range("F43").Value = -5.1721
$ws.Range("G43").Value = 44.7563
$ws.Range("H43").Value = 73.66030000000001
$ws.Range("I43").Value = -25.7509
$ws.Range("J43").Value = 17.6144
$ws.Range("K43").Value = 16.8175
$ws.Range("B45").Value = "automobiles - motorcycles / mopeds"
$ws.Range("C45").Value = -0.0317
$ws.Range("D45").Value = 0.4715
$ws.Range("E45").Value = -1.0146
$ws.Range("F45").Value = 46.6371
$ws.Range("G45").Value = 16.9168
$ws.Range("H45").Value = 86.6409
$ws.Range("I45").Value = 29.2784
$ws.Range("J45").Value = 41.8958
$ws.Range("K45").Value = 31.174
$ws.Range("C47").Value = 1.2465
$ws.Range("D47").Value = 0.1907
$ws.Range("E47").Value = -1.3275
$ws.Range("F47").Value = 22.4588
$ws.Range("C53").Value = 0.5987
$ws.Range("D53").Value = -0.0682
$ws.Range("E53").Value = 0.5395
$ws.Range("F53").Value = -33.1223
$ws.Range("C54").Value = 1.0101
$ws.Range("D54").Value = -0.3736
$ws.Range("E54").Value = 0.7557
$ws.Range("F54").Value = -49.5586
$ws.Range("B55").Value = "steel - large"
$ws.Range("C55").Value = 0.1404
$ws.Range("D55").Value = -0.5286999999999999
$ws.Range("E55").Value = -2.489
$ws.Range("F55").Value = 17.3567
$ws.Range("G55").Value = 16.1835
$ws.Range("H55").Value = 100.1913
$ws.Range("I55").Value = 119.1172
$ws.Range("J55").Value = 58.65
$ws.Range("K55").Value = 61.892
$ws.Range("C56").Value = 0.5548999999999999
$ws.Range("D56").Value = -0.5315
$ws.Range("E56").Value = -2.8092
$ws.Range("F56").Value = -0.3491
$ws.Range("B57").Value = "leather / leather products"
$ws.Range("C57").Value = 0.2516
$ws.Range("D57").Value = -0.6588000000000001
$ws.Range("E57").Value = -1.1517
$ws.Range("F57").Value = -8.3042
$ws.Range("G57").Value = -4.3507
$ws.Range("H57").Value = 43.9955
$ws.Range("I57").Value = 62.2082
$ws.Range("J57").Value = 45.2342
$ws.Range("K57").Value = 43.0938
$ws.Range("B58").Value = "computers - education"
$ws.Range("C58").Value = -1.2479
$ws.Range("D58").Value = -0.6607
$ws.Range("E58").Value = -2.5824
$ws.Range("F58").Value = -41.8669
$ws.Range("G58").Value = -30.4723
$ws.Range("H58").Value = 10.7096
$ws.Range("I58").Value = 158.5662
$ws.Range("J58").Value = 23.443
$ws.Range("K58").Value = 31.6213
$ws.Range("B59").Value = "pharmaceuticals - multinational"
$ws.Range("C59").Value = 0.0955
$ws.Range("D59").Value = -0.6828
$ws.Range("E59").Value = -4.7846
$ws.Range("F59").Value = -6.4233
$ws.Range("G59").Value = 32.1816
$ws.Range("H59").Value = 67.61750000000001
$ws.Range("I59").Value = -32.9267
$ws.Range("J59").Value = 29.9653
$ws.Range("K59").Value = 29.9724
$ws.Range("C60").Value = 0.2573
$ws.Range("D60").Value = -0.7641
$ws.Range("E60").Value = -3.0622
$ws.Range("F60").Value = -24.6728
$ws.Range("B61").Value = "printing & stationery"
$ws.Range("C61").Value = -0.1743
$ws.Range("D61").Value = -0.7969000000000001
$ws.Range("E61").Value = -3.3691
$ws.Range("F61").Value = -19.1377
$ws.Range("G61").Value = -22.033
$ws.Range("H61").Value = 57.0605
$ws.Range("I61").Value = 67.8407
$ws.Range("J61").Value = 26.8837
$ws.Range("K61").Value = 15.4547
$ws.Range("B62").Value = "entertainment / electronic media software"
$ws.Range("C62").Value = 0.0765
$ws.Range("D62").Value = -0.8054
$ws.Range("E62").Value = -0.9462
$ws.Range("F62").Value = -23.6452
$ws.Range("G62").Value = -28.2522
$ws.Range("H62").Value = 36.8911
$ws.Range("I62").Value = 96.78619999999999
$ws.Range("J62").Value = 28.0971
$ws.Range("K62").Value = 9.7845
$ws.Range("B63").Value = "pesticides / agrochemicals - indian"
$ws.Range("C63").Value = -0.0271
$ws.Range("D63").Value = -0.8717
$ws.Range("E63").Value = -8.8195
$ws.Range("F63").Value = -3.9501
$ws.Range("G63").Value = -1.0797
$ws.Range("H63").Value = -29.8151
$ws.Range("I63").Value = 31.9603
$ws.Range("J63").Value = 11.1389
$ws.Range("K63").Value = 31.7091
$ws.Range("C65").Value = 2.409
$ws.Range("D65").Value = -1.1939
$ws.Range("E65").Value = 1.9794
$ws.Range("F65").Value = -6.367
$ws.Range("C69").Value = -0.5381
$ws.Range("D69").Value = -1.7043
$ws.Range("E69").Value = -0.6508
$ws.Range("F69").Value = 55.5909
$ws.Range("C70").Value = -0.1312
$ws.Range("D70").Value = -1.8059
$ws.Range("E70").Value = -2.6223
$ws.Range("F70").Value = -36.7601
$ws.Range("B71").Value = "banks - private sector"
$ws.Range("C71").Value = 1.2008
$ws.Range("D71").Value = -1.828
$ws.Range("E71").Value = -3.2326
$ws.Range("F71").Value = -13.7051
$ws.Range("G71").Value = 21.6325
$ws.Range("H71").Value = 3.1219
$ws.Range("I71").Value = 57.8249
$ws.Range("J71").Value = 4.7572
$ws.Range("K71").Value = 4.2637
$ws.Range("B72").Value = "couriers"
$ws.Range("C72").Value = 0.3454
$ws.Range("D72").Value = -1.9449
$ws.Range("E72").Value = -7.1744
$ws.Range("F72").Value = -12.4429
$ws.Range("G72").Value = -5.3651
$ws.Range("H72").Value = -4.4588
$ws.Range("I72").Value = 54.323
$ws.Range("J72").Value = 30.7951
$ws.Range("K72").Value = 6.8212
$ws.Range("C74").Value = -0.8155
$ws.Range("D74").Value = -2.2572
$ws.Range("E74").Value = -10.0126
$ws.Range("F74").Value = -1.261
$ws.Range("C76").Value = 0.1466
$ws.Range("D76").Value = -2.6037
$ws.Range("E76").Value = -17.603
$ws.Range("F76").Value = -25.8455

$ws = $wb.Worksheets.Item("Stock List")
$ws.Range("D2").Value = 3006
$ws.Range("E2").Value = 1.4992
$ws.Range("D3").Value = 122.75
$ws.Range("E3").Value = -0.8241000000000001
$ws.Range("D4").Value = 410.85
$ws.Range("E4").Value = -1.0834
$ws.Range("D5").Value = 550.15
$ws.Range("E5").Value = 0.0364
$ws.Range("D6").Value = 1399.8
$ws.Range("E6").Value = -3.8665
$ws.Range("D12").Value = 121.94
$ws.Range("E12").Value = 0.7436
$ws.Range("D16").Value = 630.25
$ws.Range("E16").Value = 1.5059
$ws.Range("D18").Value = 1660.1
$ws.Range("E18").Value = -0.3242
$ws.Range("D19").Value = 325.5
$ws.Range("E19").Value = -0.0461
$ws.Range("D20").Value = 607.05
$ws.Range("E20").Value = -2.4819
$ws.Range("D22").Value = 142
$ws.Range("E22").Value = 5.271
$ws.Range("D23").Value = 89.79000000000001
$ws.Range("E23").Value = -0.3883
$ws.Range("D25").Value = 231.07
$ws.Range("E25").Value = -2.7647
$ws.Range("D26").Value = 217.88
$ws.Range("E26").Value = 3.565
$ws.Range("D27").Value = 438.5
$ws.Range("E27").Value = 0.712
$ws.Range("D28").Value = 316.4
$ws.Range("E28").Value = -0.4718
$ws.Range("D30").Value = 351.65
$ws.Range("E30").Value = -0.9437
$ws.Range("D31").Value = 630.2
$ws.Range("E31").Value = -1.5774
$ws.Range("D32").Value = 303.8
$ws.Range("E32").Value = 0.5294
$ws.Range("D33").Value = 656
$ws.Range("E33").Value = 6.3899
$ws.Range("D34").Value = 283.5
$ws.Range("E34").Value = 0.5141
$ws.Range("D35").Value = 979.05
$ws.Range("E35").Value = 0.6011
$ws.Range("D36").Value = 191.88
$ws.Range("E36").Value = 2.5822
$ws.Range("D38").Value = 311.25
$ws.Range("E38").Value = 4.8686
$ws.Range("D39").Value = 67.75
$ws.Range("E39").Value = 1.9564
$ws.Range("D40").Value = 339.7
$ws.Range("E40").Value = 1.5394
$ws.Range("D41").Value = 221.61
$ws.Range("E41").Value = 1.1733
$ws.Range("D43").Value = 145.68
$ws.Range("E43").Value = -0.2943
$ws.Range("D50").Value = 23.64
$ws.Range("E50").Value = 1.4157
$ws.Range("D51").Value = 105.14
$ws.Range("E51").Value = 2.4457
$ws.Range("D53").Value = 457.35
$ws.Range("E53").Value = -1.0172
$ws.Range("D56").Value = 204.9
$ws.Range("E56").Value = 1.4758
$ws.Range("D57").Value = 222.31
$ws.Range("E57").Value = -0.1392
$ws.Range("D58").Value = 317.95
$ws.Range("E58").Value = 0.5694
$ws.Range("D60").Value = 90.23999999999999
$ws.Range("E60").Value = 3.9272
$ws.Range("D61").Value = 599
$ws.Range("E61").Value = -0.7045
$ws.Range("D63").Value = 282.05
$ws.Range("E63").Value = 0.8582
$ws.Range("D68").Value = 171.06
$ws.Range("E68").Value = -0.5465
$ws.Range("D69").Value = 442.2
$ws.Range("E69").Value = 1.422
$ws.Range("D72").Value = 231.99
$ws.Range("E72").Value = 0.8300999999999999

$ws = $wb.Worksheets.Item("1 Year Return")
$ws.Range("C2").Value = 398.9201
$ws.Range("C5").Value = 274.8778
$ws.Range("C9").Value = 172.1689
$ws.Range("C10").Value = 171.0104
$ws.Range("C11").Value = 169.816
$ws.Range("C12").Value = 167.9325
$ws.Range("C13").Value = 147.8287
$ws.Range("C14").Value = 138.4828
$ws.Range("C17").Value = 127.1472
$ws.Range("C19").Value = 119.3696
$ws.Range("C20").Value = 117.1348
$ws.Range("C21").Value = 116.8576
$ws.Range("C22").Value = 108.8712
$ws.Range("C24").Value = 100.5063
$ws.Range("B25").Value = "RAMAPHO"
$ws.Range("C25").Value = 97.6888
$ws.Range("B26").Value = "VISASTEEL"
$ws.Range("C26").Value = 95.90819999999999
$ws.Range("C28").Value = 89.7954
$ws.Range("C29").Value = 88.9787
$ws.Range("C32").Value = 85.6122
$ws.Range("C33").Value = 85.1778
$ws.Range("C35").Value = 79.581
$ws.Range("C39").Value = 74.09
$ws.Range("C40").Value = 73.3296
$ws.Range("C41").Value = 73.181
$ws.Range("C43").Value = 72.0521
$ws.Range("C44").Value = 69.9337
$ws.Range("B45").Value = "ASHAPURMIN"
$ws.Range("C45").Value = 68.7901
$ws.Range("B46").Value = "GARUDA"
$ws.Range("C46").Value = 68.5898
$ws.Range("B47").Value = "SKMEGGPROD"
$ws.Range("C47").Value = 68.5164
$ws.Range("C48").Value = 68.4205
$ws.Range("B49").Value = "THANGAMAYL"
$ws.Range("C49").Value = 67.7924
$ws.Range("C50").Value = 66.1713
$ws.Range("C51").Value = 64.83320000000001
$ws.Range("C52").Value = 64.45050000000001
$ws.Range("C53").Value = 63.6475
$ws.Range("C54").Value = 62.8498
$ws.Range("C55").Value = 62.1919
$ws.Range("C56").Value = 61.6334
$ws.Range("B57").Value = "IZMO"
$ws.Range("C57").Value = 60.5278
$ws.Range("B58").Value = "IMFA"
$ws.Range("C58").Value = 60.526
$ws.Range("C59").Value = 60.2644
$ws.Range("C61").Value = 58.2054
$ws.Range("C62").Value = 57.8609
$ws.Range("C64").Value = 57.3141
$ws.Range("B65").Value = "CREDITACC"
$ws.Range("C65").Value = 57.0019
$ws.Range("B66").Value = "INDIAGLYCO"
$ws.Range("C66").Value = 56.9842
$ws.Range("C67").Value = 56.725
$ws.Range("B68").Value = "BALAJITELE"
$ws.Range("C68").Value = 56.7025
$ws.Range("B69").Value = "BSE"
$ws.Range("C69").Value = 56.5576
$ws.Range("C70").Value = 55.5909
$ws.Range("C71").Value = 55.373
$ws.Range("C72").Value = 53.6746
$ws.Range("C73").Value = 53.6017
$ws.Range("C74").Value = 53.0786
$ws.Range("C75").Value = 52.7903
